$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2.9939
$ws.Range("A3").Value = 0.89425
$ws.Range("A5").Value = 1.1513
$ws.Range("A7").Value = 1.01636
$ws.Range("A8").Value = 0.2132
$ws.Range("A9").Value = 0.14478
$ws.Range("A10").Value = 0.11326
$ws.Range("A11").Value = 0.093135
$ws.Range("A12").Value = 0.03401145000000001
$ws.Range("A13").Value = 0.01939835
$ws.Range("A14").Value = 575.3099999999999
$ws.Range("A15").Value = 12671.74
